$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-05-29"

# Update the header label in I1 ("2022 (through 05-28)" -> "2022 (through 05-29)")
$ws.Range("I1").Value = "2022 (through 05-29)"

# Update the June (row 6) 2022 total and the Total row (row 14) 2022 total
$ws.Range("I6").Value = 107
$ws.Range("I14").Value = 658
